$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "43.504.27"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +2.51%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.412.24"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +8.38%  "

$ws.Range("E4").Value = "  -0.38%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "325.75"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +11.98%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "104.78"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -5.19%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.643"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +3.32%  "

$ws.Range("E8").Value = "  -0.20%  "

$ws.Range("E9").Value = "  +7.67%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "42.30"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -3.32%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0942"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.08%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.74"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.53%  "

$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.04"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.70%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "17.26"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +15.48%  "

$ws.Range("E15").Value = "  +2.23%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "2.772.80"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +8.32%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "2.524.06"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +13.73%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "43.489.82"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +2.53%  "

$ws.Range("E19").Value = "  +4.38%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "7.39"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.01%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "75.58"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +3.63%  "

$ws.Range("E22").Value = "  +3.48%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "261.50"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +13.80%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.45"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.11%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.78"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +7.89%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "12.06"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +4.97%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.01%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "22.91"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +9.35%  "

$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "38.61"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +2.94%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "177.68"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.27%  "

$ws.Range("B31").Value = "WEMIXToken"
$ws.Range("C31").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.24"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.68%  "

$ws.Range("B32").Value = "Toncoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.17"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.77%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0938"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +6.67%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "6.00"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +6.20%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.51%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.92"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -2.74%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.0371"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -0.27%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.96"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -6.98%  "

$ws.Range("E39").Value = "  +3.30%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "2.93"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +21.13%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.61"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +23.60%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.234"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +1.07%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "124.79"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +21.30%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "69.70"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -6.15%  "

$ws.Range("E45").Value = "  +0.18%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "12.77"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.45%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "5.70"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +5.23%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.51"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +12.73%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.32"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.586.53"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +11.64%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.101"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.67%  "
